$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal_NOM_PRENOM")

# New weekly journal entry (week of 2025-12-12), filling the second block (rows 13-16)
$ws.Range("A13").Value = Get-Date -Year 2025 -Month 12 -Day 12 -Hour 0 -Minute 0 -Second 0
$ws.Range("B13").Value = "Fait la connexion entre le phidget et le PC hôte"
$ws.Range("D13").Value = 1

$ws.Range("B14").Value = "Renseignement sur les moyens de connectivité entre le serveur phidget et le frontend"
$ws.Range("D14").Value = 2.5

$ws.Range("B15").Value = "Fait le site web pour afficher les informations du robot minimalistes"
$ws.Range("D15").Value = 2.5

$ws.Range("B16").Value = "Fait le manuel d'utilisateur "
$ws.Range("D16").Value = 1

# Reflection text for this block
$ws.Range("B19").Value = "Dans la matinée, j’ai commencé par connecter le Phidget au PC hôte et vérifier que la communication fonctionnait bien. Ensuite, je me suis renseigné sur les différentes méthodes de connectivité entre le serveur Phidget et le frontend, afin de comprendre comment transmettre efficacement les données au site web.`nDans l' après-midi, j’ai créé un site web minimaliste pour afficher les informations du robot, avec une interface simple et claire. Pour finir, j’ai rédigé le manuel d’utilisateur, expliquant l’installation, la connexion au Phidget et l’utilisation du site. J'ai trouvé que ça allait bien mais on faisait pas correctement le Kanban et la documentation. J'ai été très focus à mon code et l'expérimentation"

$ws.Range("F14").Select()
